$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.816.25'
$ws.Range("D2").ClearFormats()

$ws.Range("E2").Value = '  +0.18%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.112.72'
$ws.Range("D3").ClearFormats()

$ws.Range("E3").Value = '  +0.88%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.18'
$ws.Range("D5").ClearFormats()

$ws.Range("E5").Value = '  -0.08%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.79'
$ws.Range("D6").ClearFormats()

$ws.Range("E6").Value = '  +2.23%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.109.09'
$ws.Range("D8").ClearFormats()

$ws.Range("E8").Value = '  +0.91%  '

$ws.Range("E9").Value = '  -0.33%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.46'
$ws.Range("D10").ClearFormats()

$ws.Range("E10").Value = '  -1.70%  '

$ws.Range("E11").Value = '  -0.72%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.483'
$ws.Range("D12").ClearFormats()

$ws.Range("E12").Value = '  +0.57%  '

$ws.Range("E13").Value = '  -1.44%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.01'
$ws.Range("D14").ClearFormats()

$ws.Range("E14").Value = '  +1.56%  '

$ws.Range("E15").Value = '  -1.08%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.629.85'
$ws.Range("D16").ClearFormats()

$ws.Range("E16").Value = '  +1.01%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.785.95'
$ws.Range("D17").ClearFormats()

$ws.Range("E17").Value = '  +0.12%  '

$ws.Range("E18").Value = '  -0.59%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.111.69'
$ws.Range("D19").ClearFormats()

$ws.Range("E19").Value = '  +0.93%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.26'
$ws.Range("D20").ClearFormats()

$ws.Range("E20").Value = '  +0.30%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '476.10'
$ws.Range("D21").ClearFormats()

$ws.Range("E21").Value = '  +2.03%  '

$ws.Range("E22").Value = '  -0.30%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.90'
$ws.Range("D23").ClearFormats()

$ws.Range("E23").Value = '  +4.69%  '

$ws.Range("E24").Value = '  +1.15%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.30'
$ws.Range("D25").ClearFormats()

$ws.Range("E25").Value = '  +3.50%  '

$ws.Range("E26").Value = '  -2.17%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.17'
$ws.Range("D27").ClearFormats()

$ws.Range("E27").Value = '  +0.09%  '

$ws.Range("E28").Value = '  +0.00%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.88'
$ws.Range("D29").ClearFormats()

$ws.Range("E29").Value = '  -1.42%  '

$ws.Range("E30").Value = '  -1.55%  '

$ws.Range("E31").Value = '  +0.04%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.56'
$ws.Range("D32").ClearFormats()

$ws.Range("E32").Value = '  +1.34%  '

$ws.Range("E33").Value = '  +0.91%  '

$ws.Range("E34").Value = '  -7.95%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").ClearFormats()

$ws.Range("E35").Value = '  +0.08%  '

$ws.Range("E36").Value = '  -0.49%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.972'
$ws.Range("D37").ClearFormats()

$ws.Range("E37").Value = '  -2.93%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '47.17'
$ws.Range("D38").ClearFormats()

$ws.Range("E38").Value = '  +0.25%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '50.18'
$ws.Range("D39").ClearFormats()

$ws.Range("E39").Value = '  -0.16%  '

$ws.Range("E40").Value = '  -3.18%  '

$ws.Range("E41").Value = '  -2.47%  '

$ws.Range("E42").Value = '  -0.02%  '

$ws.Range("E43").Value = '  +0.05%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '387.87'
$ws.Range("D44").ClearFormats()

$ws.Range("E44").Value = '  +1.25%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.833.44'
$ws.Range("D45").ClearFormats()

$ws.Range("E45").Value = '  +2.83%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0355'
$ws.Range("D46").ClearFormats()

$ws.Range("E46").Value = '  -1.49%  '

$ws.Range("E47").Value = '  -8.89%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '135.32'
$ws.Range("D48").ClearFormats()

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.73'
$ws.Range("D50").ClearFormats()

$ws.Range("E50").Value = '  +0.40%  '

$ws.Range("E51").Value = '  -1.93%  '
